# Auto-generated Excel COM-interop edit script
# Updates the cryptos list (Price / Volume(1h) columns, and a row swap)
# to match the Thu Nov  2 07:10:31 UTC 2023 GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing a Text number format so
# numeric-looking strings (e.g. "1.00", "229.50") keep their exact original
# text representation instead of Excel normalising them into a Double. The
# temporary "@" format is removed again immediately afterwards (restoring the
# cell to the plain "Normal" style) so formatting is left untouched.
function Set-TextCell($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "35.390.43"
$ws.Range("E2").Value = "  +2.74%  "

# Row 3
Set-TextCell "D3" "1.842.10"
$ws.Range("E3").Value = "  +2.03%  "

# Row 4
$ws.Range("E4").Value = "  +0.24%  "

# Row 5
Set-TextCell "D5" "229.50"
$ws.Range("E5").Value = "  +2.14%  "

# Row 6
$ws.Range("E6").Value = "  +2.72%  "

# Row 7
$ws.Range("E7").Value = "  +0.21%  "

# Row 8
Set-TextCell "D8" "43.31"
$ws.Range("E8").Value = "  +13.36%  "

# Row 9
Set-TextCell "D9" "0.308"
$ws.Range("E9").Value = "  +7.14%  "

# Row 10
Set-TextCell "D10" "0.0695"
$ws.Range("E10").Value = "  +3.56%  "

# Row 11
$ws.Range("E11").Value = "  +3.63%  "

# Row 12
Set-TextCell "D12" "2.107.78"
$ws.Range("E12").Value = "  +2.01%  "

# Row 13
Set-TextCell "D13" "1.843.03"
$ws.Range("E13").Value = "  +2.17%  "

# Row 14
Set-TextCell "D14" "11.29"
$ws.Range("E14").Value = "  +2.29%  "

# Row 15
Set-TextCell "D15" "0.673"
$ws.Range("E15").Value = "  +7.29%  "

# Row 16
Set-TextCell "D16" "4.67"
$ws.Range("E16").Value = "  +6.65%  "

# Row 17
Set-TextCell "D17" "35.383.56"
$ws.Range("E17").Value = "  +2.80%  "

# Row 18
Set-TextCell "D18" "70.16"
$ws.Range("E18").Value = "  +3.30%  "

# Row 19
$ws.Range("E19").Value = "  +3.45%  "

# Row 20
Set-TextCell "D20" "244.43"
$ws.Range("E20").Value = "  +1.21%  "

# Row 21
Set-TextCell "D21" "12.10"
$ws.Range("E21").Value = "  +9.42%  "

# Row 22
Set-TextCell "D22" "4.66"
$ws.Range("E22").Value = "  +14.02%  "

# Row 23
$ws.Range("E23").Value = "  +0.28%  "

# Row 24
$ws.Range("E24").Value = "  +0.36%  "

# Row 25
Set-TextCell "D25" "168.87"
$ws.Range("E25").Value = "  -1.36%  "

# Row 26
$ws.Range("E26").Value = "  +2.75%  "

# Row 27
Set-TextCell "D27" "17.77"
$ws.Range("E27").Value = "  +2.24%  "

# Row 28
$ws.Range("E28").Value = "  +1.86%  "

# Row 30
$ws.Range("E30").Value = "  +0.16%  "

# Row 31
Set-TextCell "D31" "3.245.32"
$ws.Range("E31").Value = "  +33.57%  "

# Row 32
$ws.Range("E32").Value = "  +6.50%  "

# Row 33
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D33" "4.07"
$ws.Range("E33").Value = "  +6.04%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D34" "3.92"
$ws.Range("E34").Value = "  +4.43%  "

# Row 35
$ws.Range("E35").Value = "  +2.64%  "

# Row 36
Set-TextCell "D36" "94.88"
$ws.Range("E36").Value = "  +14.49%  "

# Row 37
Set-TextCell "D37" "0.688"
$ws.Range("E37").Value = "  +7.64%  "

# Row 38
Set-TextCell "D38" "1.346.13"
$ws.Range("E38").Value = "  +2.17%  "

# Row 39
Set-TextCell "D39" "1.09"
$ws.Range("E39").Value = "  +3.08%  "

# Row 40
$ws.Range("E40").Value = "  +5.90%  "

# Row 41
$ws.Range("E41").Value = "  +3.69%  "

# Row 42
Set-TextCell "D42" "1.00"
$ws.Range("E42").Value = "  +6.26%  "

# Row 43
$ws.Range("E43").Value = "  +3.75%  "

# Row 44
Set-TextCell "D44" "14.78"
$ws.Range("E44").Value = "  +7.89%  "

# Row 45
$ws.Range("E45").Value = "  +0.62%  "

# Row 46
$ws.Range("E46").Value = "  -0.03%  "

# Row 47
Set-TextCell "D47" "6.25"
$ws.Range("E47").Value = "  +8.69%  "

# Row 48
$ws.Range("E48").Value = "  +1.56%  "

# Row 49
Set-TextCell "D49" "2.009.18"
$ws.Range("E49").Value = "  +2.15%  "

# Row 50
$ws.Range("E50").Value = "  +0.24%  "

# Row 51
Set-TextCell "D51" "103.00"
$ws.Range("E51").Value = "  +0.86%  "

